# Fix Training Data Issue (#48)
# Data was taken from 1 day off due to way NBA stats were shown.
# Column BF ("Date") held the sheet's own file-name-derived label
# "5-15-2013-14" for every data row; correct it to the real game date
# "2014-05-15" (ISO yyyy-mm-dd) for rows 2-31.
#
# NOTE: assigning a plain ISO-looking string via .Value/.Value2 makes
# Excel auto-convert it to a date serial (and pulls in a new number-format
# style). To keep it as literal text with no style change, write it as a
# formula that evaluates to the literal string, then convert the range to
# values in place via copy / paste-special-values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$firstRow = 2
$lastRow = 31
$oldDate = "5-15-2013-14"
$newDate = "2014-05-15"

for ($row = $firstRow; $row -le $lastRow; $row++) {
    $cell = $ws.Range("BF$row")
    if ($cell.Value2 -eq $oldDate) {
        $cell.Formula = '="' + $newDate + '"'
    }
}

$rng = $ws.Range("BF$firstRow`:BF$lastRow")
$rng.Copy()
$rng.PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$excel.CutCopyMode = $false
